$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "288.50"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-9.69%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "40.43"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-2.26%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.029"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07288"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-5.92%"
$ws.Range("B6").Value = "FTXToken"
$ws.Range("C6").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.522"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-12.84%"
$ws.Range("B7").Value = "MXToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9174"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-2.85%"
$ws.Range("B8").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C8").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.1177"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-4.98%"
$ws.Range("B9").Value = "WazirX"
$ws.Range("C9").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1725"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-8.90%"
$ws.Range("B10").Value = "MandalaExchangeToken"
$ws.Range("C10").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08650"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-5.55%"
$ws.Range("B11").Value = "BitrueCoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.04177"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-2.79%"
$ws.Range("B12").Value = "BitMartToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.1053"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "0.12%"
$ws.Range("B13").Value = "BitForexToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.001267"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-1.73%"
$ws.Range("B14").Value = "TigerCash"
$ws.Range("C14").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.005779"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-2.11%"
$ws.Range("B15").Value = "LEO"
$ws.Range("C15").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.390"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "1.48%"
$ws.Range("B16").Value = "GateToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.284"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-1.34%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.3318"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-1.20%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.883"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "1.33%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1343"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-0.82%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.2887"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "2.27%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.03874"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-3.83%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.001270"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "0.13%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.003854"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-6.58%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0001282"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "0.97%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0003729"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02309"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "-9.86%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.04954"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-7.11%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.006743"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "238.35%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007676"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-1.02%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1273"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-3.27%"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "4.41%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.007070"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-14.28%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.3113"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-1.73%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006439"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-3.64%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.18%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.03500"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-82.62%"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "0.03%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002104"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.18%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002004"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.18%"
